$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the car name / filename values in column A
$ws.Range("A2").Value = "ford-tourneo-custom-2024-1019-safety pack.xlsx"
$ws.Range("A3").Value = "zeekr-001-2024-1037.xlsx"
$ws.Range("A4").Value = "maxus-mifa-7-2024-1060.xlsx"
$ws.Range("A5").Value = "ford-tourneo-custom-2024-1019-standard.xlsx"
$ws.Range("A6").Value = "vw-passat-2024-1021.xlsx"
$ws.Range("A7").Value = "skoda-kodiaq-2024-1074.xlsx"
$ws.Range("A8").Value = "bmw-x2-2022-1065.xlsx"
$ws.Range("A9").Value = "renault-rafale-hev-2022-1073.xlsx"
$ws.Range("A10").Value = "mercedes-benz-e-class-2024-1064.xlsx"
$ws.Range("A11").Value = "suzuki-swift-2024-1036.xlsx"
$ws.Range("A12").Value = "dacia-duster-2024-1075.xlsx"
$ws.Range("A13").Value = "renault-espace-2022-1072.xlsx"
$ws.Range("A14").Value = "Toyota C-HR-2024-1042.xlsx"
$ws.Range("A15").Value = "Honda CR-V-2024-1050-Safety Pack.xlsx"
$ws.Range("A16").Value = "NIO EL6-2024-1054.xlsx"
$ws.Range("A17").Value = "Honda CR-V-2024-1050-Standard.xlsx"

# Widen column A to fit the longer filenames (results in stored width of 48)
$ws.Columns.Item(1).ColumnWidth = 47.2
